# Modify the project summary website row to show a link if there is a
# website - adds three new portfolio project rows (MechaCar Statistics,
# Bike Sharing, Employee Database) plus a new "website_name" column that
# is populated (with a hyperlink) only when a project has a website.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J: website_name header -------------------------------
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "website_name"

# --- Row 13: MechaCar Statistics (order 11) ---------------------------
$ws.Range("A12:H12").Copy()
$ws.Range("A13:H13").PasteSpecial(-4122)

$ws.Range("A13").Value = "mechacar"
$ws.Range("B13").Value = 11
$ws.Range("C13").Value = "r"
$ws.Range("D13").Value = "MechaCar Statistics"
$ws.Range("E13").Value = "https://github.com/cdpeters/MechaCar-statistical-analysis-R"
$ws.Range("F13").Value = "R"
$ws.Range("G13").Value = "dplyr"
$ws.Range("H13").Value = "linear regression, t-test, hypothesis testing, study design"

# --- Row 14: Bike Sharing (order 12) -----------------------------------
$ws.Range("A12:G12").Copy()
$ws.Range("A14:G14").PasteSpecial(-4122)
$ws.Range("E12").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("A12").Copy()
$ws.Range("J14").PasteSpecial(-4122)

$ws.Range("A14").Value = "bike_sharing"
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = "tableau"
$ws.Range("D14").Value = "Bike Sharing"
$ws.Range("E14").Value = "https://github.com/cdpeters/bike-sharing-tableau"
$ws.Range("F14").Value = "Python"
$ws.Range("G14").Value = "pandas, tableau"
$ws.Range("I14").Value = "https://public.tableau.com/views/NYC_CitiBike_Challenge_16506220556720/August2019NYCCitibikeStudy?:language=en-US&:display_count=n&:origin=viz_share_link"
$ws.Range("J14").Value = "NYC CitiBike Challenge"

# --- Row 15: Employee Database (order 16) ------------------------------
$ws.Range("A12:H12").Copy()
$ws.Range("A15:H15").PasteSpecial(-4122)

$ws.Range("A15").Value = "employee_db"
$ws.Range("B15").Value = 16
$ws.Range("C15").Value = "sql"
$ws.Range("D15").Value = "Employee Database"
$ws.Range("E15").Value = "https://github.com/cdpeters/employee-database-postgresql"
$ws.Range("F15").Value = "SQL"
$ws.Range("G15").Value = "postgresql"
$ws.Range("H15").Value = "database"

# --- Hyperlinks (added in the same order they appear in the file) -----
$ws.Hyperlinks.Add($ws.Range("E13"), "https://github.com/cdpeters/MechaCar-statistical-analysis-R")
$ws.Hyperlinks.Add($ws.Range("E14"), "https://github.com/cdpeters/bike-sharing-tableau")
$ws.Hyperlinks.Add($ws.Range("I14"), "https://public.tableau.com/views/NYC_CitiBike_Challenge_16506220556720/August2019NYCCitibikeStudy?:language=en-US&:display_count=n&:origin=viz_share_link")
$ws.Hyperlinks.Add($ws.Range("E15"), "https://github.com/cdpeters/employee-database-postgresql")

# --- Column J width (bestFit-like autosize for the new column) --------
$ws.Columns.Item(10).ColumnWidth = 23.05

# --- Selection ends on the newly added row's repo cell -----------------
$ws.Range("E15").Select()
